$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data (now B:F) right
$ws.Columns.Item(1).Insert()

# Copy the header style/format from B1 (the old A1) onto the new A1
$ws.Cells.Item(1, 2).Copy()
$ws.Cells.Item(1, 1).PasteSpecial(-4122)

# Set header text for the new ID column
$ws.Cells.Item(1, 1).Value = "ID"

# Fill in ID values for rows 2-25
$ids = @("Hb 2", "Hb 3", "S 24", "S 28", "Hb 107", "Hb 66", "Hb 69", "Hb 95", "Hb 99", "Hb 92", "Hb 40", "Hb 41", "S 11", "Hb 57", "S 21", "S 22", "S 3", "S 4", "S 5", "Hb 74", "Hb 79", "Hb 32", "S 15", "S 16")

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $ids[$i]
}
